# "Fixed data numbers for 0.15.9"
#
# The underlying data table (Output / Input / Cost) gets a handful of
# corrected Cost values, then the whole table is re-sorted by the Output
# column (ascending, with header) and an AutoFilter is turned on over the
# table range - exactly what Excel records when you select the range,
# fix some numbers, and use Data > Sort / Data > Filter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the mis-entered Cost values (by their original row) ---
$ws.Cells.Item(32, 3).Value = 5    # Piercing Rounds Magazine / Copper Plate: 4 -> 5
$ws.Cells.Item(34, 3).Value = 4    # Firearm Magazine / Iron Plate: 2 -> 4
$ws.Cells.Item(51, 3).Value = 15   # Electric Engine Unit / Lubricant: 2 -> 15
$ws.Cells.Item(52, 3).Value = 10   # Electric Furnace / Steel: 15 -> 10
$ws.Cells.Item(56, 3).Value = 30   # Yellow Science Pack / Copper Cable: 40 -> 30
$ws.Cells.Item(58, 3).Value = 3    # Yellow Science Pack / Blue Circuit: 1 -> 3
$ws.Cells.Item(62, 3).Value = 20   # Battery / Sulfuric Acid: 2 -> 20
$ws.Cells.Item(65, 3).Value = 100  # Sulfuric Acid / Water: 10 -> 100
$ws.Cells.Item(66, 3).Value = 5    # Blue Circuit / Sulfuric Acid: 0.5 -> 5

# --- 2. Turn on filtering for the table range ---
$tableRange = $ws.Range("A1:C75")
$tableRange.AutoFilter() | Out-Null

# Excel records the filtered range as a hidden sheet-local defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$75")
$filterName.Visible = $false

# --- 3. Re-sort the whole table by Output (column A), ascending, header kept ---
$ws.Sort.SortFields.Clear() | Out-Null
$ws.Sort.SortFields.Add($ws.Range("A1:A75")) | Out-Null
$ws.Sort.SetRange($tableRange) | Out-Null
$ws.Sort.Header = 1
$ws.Sort.Apply() | Out-Null

# Leave the selection somewhere sensible on the sheet, as Excel would.
$ws.Range("C73").Select() | Out-Null
